# Mixorcerers TODO workbook update
# - Adds a new "Pitch" sheet (the game design pitch / GDD) after "Backlog"
# - Makes "Pitch" the active tab
# - Tidies up the "Backlog" sheet (removes the stray blank row 3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Backlog sheet cleanup: row 3 only ever held a blank, default-styled
#    A3 cell - drop it so the sheet data goes straight from row 2 to row 4.
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("A3").Clear()

# ---------------------------------------------------------------------------
# 2. Add the "Pitch" sheet right after "Backlog" and make it active.
# ---------------------------------------------------------------------------
$pitch = $wb.Worksheets.Add($null, $backlog)
$pitch.Name = "Pitch"

# Column widths (~34.67 and ~115.2 characters).
$pitch.Columns.Item(1).ColumnWidth = 33.8
$pitch.Columns.Item(2).ColumnWidth = 114.4

# ---------------------------------------------------------------------------
# 3. Content - column A labels / column B answers, grouped in three blocks
#    separated by blank rows (6 and 12).
# ---------------------------------------------------------------------------

# --- Block 1: rows 1-5 (working title / concept / genre / audience / USP)
$pitch.Range("A1").Value = "Working title"
$pitch.Range("B1").Value = "Mixorcerers"

$pitch.Range("A2").Value = "Concept statement"
$pitch.Range("B2").Value = "Mixorcerers is a turn based 1v1 game that combines rts elements with the unlimited power of a mage. "

$pitch.Range("A3").Value = "Genre(s)"
$pitch.Range("B3").Value = "Strategy"

$pitch.Range("A4").Value = "Target audience"
$pitch.Range("B4").Value = "Universal ESRB " + [char]0x2013 + " suitable 13 and up"

$pitch.Range("A5").Value = "Unique Selling Points"
$pitch.Range("B5").Value = "Unique Power Fantasy, Interesting Combat System"

# --- Block 2: rows 7-11 (player experience / visuals / world / monetization / platform)
$pitch.Range("A7").Value = "Player Experience and Game POV"
$pitch.Range("B7").Value = "Player is a mage. Setting is cheerful medieaval. Archmage fantasy. Power, scheming, accomplishment, gumption"

$pitch.Range("A8").Value = "Visual and Audio Style"
$pitch.Range("B8").Value = "Reminiscent of FE7"

$pitch.Range("A9").Value = "Game World Fiction"
$pitch.Range("B9").Value = "1v1 arena, maybe some mcguffin, some kinda astral projection to keep the constant battles canon"

$pitch.Range("A10").Value = "Monetization"
$pitch.Range("B10").Value = "Feed me"

$pitch.Range("A11").Value = "Platform(s), Technology, and Scope (brief)"
$pitch.Range("B11").Value = "PC and mobile. Godot. Maybe a couple more months, team of me. Couple more months? "

# --- Block 3: rows 13-16 (core loops / objectives / systems / interactivity)
$pitch.Range("A13").Value = "Core Loops"
$pitch.Range("B13").Value = "1v1 laddering is a classic loop of self improvement. The rich variety of strategies will also cause people to come back again. The different types of maps and some kind of randomness in gaining orbs will stop the meta from getting stale. Currently im thinking of bases to capture as well as worker harassment for magycke. I expect to see a dedicated but small community that will play the game for many years, and hopefully I can hand off balance and such to them eventually"

$pitch.Range("A14").Value = "Objectives and Progression"
$pitch.Range("B14").Value = "Theres a tutorial, and then the ladder. There will be a leaderboard.While battling people may find stuff that improves their lobby screen which is like their personal den.In the future I may add puzzles if Im jobless. "

$pitch.Range("A15").Value = "Game Systems"
$pitch.Range("B15").Value = "Multiplayer Server,Replays, Client " + [char]0x2013 + " Game, Menus, DisplayShowcase , Main webserver that does leaderboard, display showcase,matchmaking"

$pitch.Range("A16").Value = "Interactivity"
$pitch.Range("B16").Value = "The interactivity is nice keyboard shortcuts for everything. The player moves the cursor which interacts with the units and structures, then uses the menu to execute commands. The combat works by casting spells on the map"

# ---------------------------------------------------------------------------
# 4. Row heights - rows with wrapped / multi-line content are taller.
# ---------------------------------------------------------------------------
$pitch.Rows.Item(2).RowHeight = 23.9
$pitch.Rows.Item(13).RowHeight = 49.55
$pitch.Rows.Item(14).RowHeight = 26.45
$pitch.Rows.Item(16).RowHeight = 28.6

# ---------------------------------------------------------------------------
# 5. Word wrap on the long-form answers.
# ---------------------------------------------------------------------------
$pitch.Range("B2").WrapText = $true
$pitch.Range("B13").WrapText = $true
$pitch.Range("B14").WrapText = $true
$pitch.Range("B16").WrapText = $true

# ---------------------------------------------------------------------------
# 6. Fill colors per block.
#    Block 1 (rows 1-5): fg FFFF6D / bg FFFFCC
#    Block 2 (rows 7-11): fg DDE8CB / bg F7D1D5
#    Block 3 (rows 13-16): fg F7D1D5 / bg DDE8CB
# ---------------------------------------------------------------------------
$block1 = $pitch.Range("A1:B5")
$block1.Interior.Color = 7208959
$block1.Interior.PatternColor = 13434879

$block2 = $pitch.Range("A7:B11")
$block2.Interior.Color = 13363421
$block2.Interior.PatternColor = 14012919

$block3 = $pitch.Range("A13:B16")
$block3.Interior.Color = 14012919
$block3.Interior.PatternColor = 13363421

Write-Output "Pitch sheet added"
